$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$b64 = "aDE6IE1pZ2h0IG9mIEFsbGFoCnA6IFRvZGF5IFF1cmFuIHRhbGtzIHByZWNpc2VseSBhYm91dCB0aGUgZGlzdHJpYnV0aW9uIG9mIGluaGVyaXRhbmNlIHRvIG91ciBjbG9zZSBvbmVzLiBJZiBDaGlsZCBkaWVzLCBtb3RoZXIgZ2V0cyAxLzZ0aCBvZiBoZXIgc2hhcmUuIElmIG1vdGhlciBpcyBhIHdpZG93ZXIsIHNoZSBnZXRzIDEvOHRoIG9mIGhlciBzaGFyZS4gQWZ0ZXIgZWFjaCBkZWNpc2lvbiwgUXVyYW4gc3BlYWtzIG9mIHRoZSBtaWdodCBvZiB0aGUgQ3JlYXRvci4gRm9yIGV4YW1wbGUKcC5iLWxlZnQ6IDxiPjEuIFZlcnNlIDQ8L2I+LiBBbmQgc3VmZmljaWVudCBpcyBBbGxhaCBhcyBBY2NvdW50YW50LgpwLmItbGVmdDogPGI+Mi4gVmVyc2UgOTwvYj4uIFNvIGxldCB0aGVtIGZlYXIgQWxsYWggYW5kIHNwZWFrIHdvcmRzIG9mIGFwcHJvcHJpYXRlIGp1c3RpY2UuCnAuYi1sZWZ0OiA8Yj4zLiBWZXJzZSAxMDwvYj4uIEFuZCB0aGV5IHdpbGwgYmUgYnVybmVkIGluIGEgQmxhemUuCnAuYi1sZWZ0OiA8Yj40LiBWZXJzZSAxMTwvYj4uIEluZGVlZCwgQWxsYWggaXMgZXZlciBLbm93aW5nIGFuZCBXaXNlLgpwLmItbGVmdDogPGI+NS4gVmVyc2UgMTI8L2I+LiBBbGxhaCBpcyBLbm93aW5nIGFuZCBGb3JiZWFyaW5nLgpwOiBUaGUgcmh5dGhtIG9mIHRoZXNlIHNlbnRlbmNlcyBpcyBmb3VuZCBpbiBlbnRpcmUgUXVyYW4uIFRoZXJlIGFyZSBkZWNpc2lvbnMsIGluc3RydWN0aW9ucywgb2JsaWdhdGlvbnMgYW5kIGltbWVkaWF0ZWx5IGFmdGVyIHRoZXNlIHJ1bGluZ3MsIHRoZXJlIGlzIHRoZSBtaWdodCBvZiBBbGxhaC4gCmgzOiBCZWxpZXZpbmcgaW4gUXVyYW4KcDogRmlyc3QgbWlyYWNsZSBpbiBhbnlib2R54oCZcyBsaWZlIGlzLCBtYXkgYmUsIGJpcnRoLiBZb3UgY29tZSBvdXQgb2YgeW91ciA8Yj5tb3RoZXLigJlzIGJlbGx5PC9iPi4gWW91IGdyb3cgYXMgeW91IGVhdDxiPiBmb29kPC9iPi4gQWZ0ZXIgc3BlbmRpbmcgdGhlIHJlc3Qgb2YgdGhlIGxpZmUgaW4gdGhpcyB3b3JsZCwgeW91PGI+IGRpZTwvYj4uCnA6IER1cmluZyBlYWNoIG1pcmFjbGUsIHNjaWVuY2Ugc3RyaXZlcyB0byBleHBsYWluIGl0LiBUaGUgYmVsbHkgZ2l2ZXMgdGhlIGNoaWxkLCB3ZSBzdGFydCBsYWJlbGxpbmcgaXQgYSBuYXR1cmFsIGFjdC4gU2NpZW5jZSBsYWJlbHMgaXQgc2VxdWVuY2Ugb2YgbWlyYWN1bG91cyBldmVudHMuIApwOiBXZSBsYWJlbCBmb29kLCBuYXR1cmFsIGluZ3JlZGllbnRzLCB0byBydW4gaHVtYW4gYm9keS4gV2UgZ3JvdywgbGFiZWxsaW5nIGl0IG5hdHVyYWwuIFdlIGRpZSBhbmQgbGFiZWwgZGVhdGggYSByZXN1bHQgb2YgTmF0dXJhbCBDYXVzZXMuCnA6IFRoZSBlYXJ0aCBpcyBhbiBleHRyYSBvcmRpbmFyeSBwaWVjZSBvZiBuYXR1cmUuIEVhY2ggcGFzc2luZyBzZWNvbmQsIHRoZSBlYXJ0aCBpcyBzdXJyb3VuZGVkIHdpdGggYW5vdGhlciBtaXJhY2xlLiBXZSB3aXRuZXNzIHRoaXMgbWlyYWNsZSBhbmQgPGI+d2UgcnVuIHRvIHByb3ZlIGl0IGlzIG5hdHVyYWwgb3Igc2NpZW50aWZpYzwvYj4uIApwOiBXZSBzdHJpdmUgaGFyZCB0byBtYWtlIG91cnNlbHZlcyBiZWxpZXZlIGl0IGlzIGp1c3QgYW5vdGhlciBhY3Qgb2YgbmF0dXJlLiAKaDM6IE5hdHVyZSA9ID8KcDogTmF0dXJlIGlzIG5vdCBhcnRpZmljaWFsLiBJdCBpcyBhbnl0aGluZyBicmVhdGhpbmcuIFRoZSBlYXJ0aCwgaXMgYnJlYXRoaW5nLiBBbmltYWxzIGdyYXppbmcgb24gaXQgYXJlIGFsc28gYnJlYXRoaW5nLiBXZSBiZWxpZXZlIHRoZSBuYXR1cmUgb2YgZWFydGggaXMgY2F1c2luZyB0aGlzIHBsYW5ldCB0byByb3RhdGUuIFRoZXJlIGlzIGdyYXZpdHkgbGlua2VkIHdpdGggdGhlIGVhcnRoLCB0aGF0IGtlZXBzIHRoZSBhbmltYWxzIGZyb20gZmFsbGluZyBvZmYgaXQuIEFueXRoaW5nIHRoYXQgY2FuIG5vdCBiZSByZWNvbnN0cnVjdGVkIGlzIG5hdHVyZS4KcDogV2UgbGFiZWwgbmF0dXJlIGFzIGEgcmVzdWx0IG9mIGNvaW5jaWRlbmNlIG9mIG1hdHRlciBhbmQgbW9sZWN1bGUuCmgzOiBEaWQgYW55b25lIGZpbmQgdGhlIGFuc3dlcj8KcDogRXZlcnkgY2VsbHVsYXIgc3RydWN0dXJlIGluc2lkZSBuYXR1cmUgaXMgYSBjb21iaW5hdGlvbiBvZiBlbGVjdHJvbnMsIHByb3RvbnMgYW5kIG5ldXRyb25zLiBUaGVzZSBhcmUgbWludXRlIHBhcnRpY2xlcyB0aGF0IG1ha2UgZXZlcnkgbW9sZWN1bGUuIFNvbWUgZ3JvdyB0byBiZWNvbWUgYSBzcGVybSBhbmQgZmV3IGdyb3cgdG8gYmVjb21lIGVhcnRoLiBGcm9tIG1pbnV0ZSB0byBiaWdnZXN0IHN0cnVjdHVyZSBpcyBhIGNvbWJpbmF0aW9uIG9mIGNlbGxzLiAKcDogTmV4dCBxdWVzdGlvbnMgdGhhdCBuZWVkIHRvIGJlIGFza2VkCnAuYi1sZWZ0OiAxLiBXaHkgc29tZSBlbGVjdHJvbnMgZ3JvdyB0byBsb29rIGxpa2Ugc3Blcm0/CnAuYi1sZWZ0OiAyLiBIb3cgc3Blcm0gZGVjaWRlcyBpdCBnb2VzIHRvIHdvbWI/CnAuYi1sZWZ0OiAzLiBIb3cgZG9lcyBhIGJhYnkgc3RheXMgYWxpdmUgaW5zaWRlIGEgY2xvc2VkIGNoYW1iZXI/CnAuYi1sZWZ0OiA0LiBIb3cgZG9lcyBhIHNwZXJtIGdyb3dzIHRvIGEgaHVtYW4gYmVpbmc/CnAuYi1sZWZ0OiA1LiBXaGF0IGtlZXBzIHRoZSBodW1hbiBmcm9tIGZhbGxpbmcgb2ZmIHRoZSBwbGFuZXQgRWFydGg/CnAuYi1sZWZ0OiA2LiBIb3cgbWFueSBtb3JlIEVhcnRo4oCZcyBleGlzdCBpbiBhIHVuaXZlcnNlIG9mIGNvaW5jaWRlbmNlcz8KcC5iLWxlZnQ6IDcuIEhvdyBiaWcgaXMgdW5pdmVyc2U/CmgzOiBEaWQgd2UgZmluZCB0aGUgYW5zd2VyIHRvIG91ciBxdWVzdGlvbnM/CnA6IE5vLiBPdXIgbm8gcXVlc3Rpb24gYWJvdXQgbmF0dXJlIGdvZXMgYW5zd2VyZWQuIEVhY2ggcXVlc3Rpb24gdGFrZXMgdXMgZnVydGhlciBkb3duIGEgc2VyaWVzIG9mIHJhYmJpdCBob2xlcyBhbmQgdGhlc2UgcXVlc3Rpb25zIG5ldmVyIHN0b3AuICAKcDogV2UgZGlkIG5vdCBmaW5kIG91dCB0aGUgcmVhbGl0eSBiZWhpbmQgYSBzcGVybS4gV2h5IGl0IGVudGVycyB0aGUgd29tYj8gV2hhdCBtYWtlcyBhIGJpcnRoPyBFYWNoIHF1ZXN0aW9uIGlzIG9ubHkgYW4gYXJndW1lbnQuCnA6IFRoZSBhbnN3ZXJzIHdlIGZpbmQgYXJlIGxpdGVyYWwgc2NpZW50aWZpYyBzdG9yaWVzLiBTY2llbmNlIG5ldmVyIHN1cnJlbmRlcnMsIGl0IGlzIGEgbmFtZSBvZiByYWJiaXQgaG9sZS4gCmgzOiBDb25jbHVzaW9uCnA6IFdlIGRvIG5vdCBhc2sgMSBxdWVzdGlvbi4gCnAuYi1sZWZ0OiA8Yj5XaGF0IGlmIHdlIGRpZSBhbmQgdGhlcmUgaXMgYSBHb2Q/PC9iPgpwOiBXZSBkbyBub3Qgc3VibWl0IHRvIHRoZSB2YWN1dW0gY3JlYXRlZCBieSBvdXIgdW5hbnN3ZXJlZCBxdWVzdGlvbnMuIEV4ZWN1dGlvbiBvbiBRdXJhbmljIGluc3RydWN0aW9ucyBpcyBvbmx5IHBvc3NpYmxlIGFmdGVyIHdlIGZpbmQgQWxsYWguIApwOiBPbmNlIGZvdW5kOyBzdWJtaXQsIHN1cnJlbmRlciBhbmQgbW92ZSBuZWFyIEhpbS4gSGUgaXMgZ29pbmcgdG8gbWVldCB1cyB0aGUgZGF5IHdlIGRpZSBpbiBmZXcgZGF5cyBvciBtb250aHMuIEEgTWlnaHQgaXMgZXhwZWN0aW5nIHVzLiA8Yj5JdCBjYW4gaGFwcGVuIHRoYXQgdGhlIGRheSBvZiBvdXIgZGVhdGggYmVjb21lcyB0aGUg4oCYSGFwcGllc3QgZGF54oCZIG9mIG91ciBsaXZlcy4gPC9iPgpwLm5vdGU6IEd1eXMsIGlmIHlvdSBsaWtlIHRoaXMgcHJvamVjdC4gUGxlYXNlIGZvbGxvdyB0aGlzIHByb2plY3QncyBwYWdlIG9uIHR3aXR0ZXIuIDxhIGhyZWY9Imh0dHBzOi8vdHdpdHRlci5jb20vemFrYXRsaXN0cyI+Q2xpY2sgaGVyZSB0byBnbyB0byB0aGUgdHdpdHRlciBwYWdlPC9hPi4gU3Vic2NyaWJlIGJlbG93IHRvIHJlY2VpdmUgdXBkYXRlcy5oMTogTWlnaHQgb2YgQWxsYWgKcDogVG9kYXkgUXVyYW4gdGFsa3MgcHJlY2lzZWx5IGFib3V0IHRoZSBkaXN0cmlidXRpb24gb2YgaW5oZXJpdGFuY2UgdG8gb3VyIGNsb3NlIG9uZXMuIElmIENoaWxkIGRpZXMsIG1vdGhlciBnZXRzIDEvNnRoIG9mIGhlciBzaGFyZS4gSWYgbW90aGVyIGlzIGEgd2lkb3dlciwgc2hlIGdldHMgMS84dGggb2YgaGVyIHNoYXJlLiBBZnRlciBlYWNoIGRlY2lzaW9uLCBRdXJhbiBzcGVha3Mgb2YgdGhlIG1pZ2h0IG9mIHRoZSBDcmVhdG9yLiBGb3IgZXhhbXBsZQpwLmItbGVmdDogPGI+MS4gVmVyc2UgNDwvYj4uIEFuZCBzdWZmaWNpZW50IGlzIEFsbGFoIGFzIEFjY291bnRhbnQuCnAuYi1sZWZ0OiA8Yj4yLiBWZXJzZSA5PC9iPi4gU28gbGV0IHRoZW0gZmVhciBBbGxhaCBhbmQgc3BlYWsgd29yZHMgb2YgYXBwcm9wcmlhdGUganVzdGljZS4KcC5iLWxlZnQ6IDxiPjMuIFZlcnNlIDEwPC9iPi4gQW5kIHRoZXkgd2lsbCBiZSBidXJuZWQgaW4gYSBCbGF6ZS4KcC5iLWxlZnQ6IDxiPjQuIFZlcnNlIDExPC9iPi4gSW5kZWVkLCBBbGxhaCBpcyBldmVyIEtub3dpbmcgYW5kIFdpc2UuCnAuYi1sZWZ0OiA8Yj41LiBWZXJzZSAxMjwvYj4uIEFsbGFoIGlzIEtub3dpbmcgYW5kIEZvcmJlYXJpbmcuCnA6IFRoZSByaHl0aG0gb2YgdGhlc2Ugc2VudGVuY2VzIGlzIGZvdW5kIGluIGVudGlyZSBRdXJhbi4gVGhlcmUgYXJlIGRlY2lzaW9ucywgaW5zdHJ1Y3Rpb25zLCBvYmxpZ2F0aW9ucyBhbmQgaW1tZWRpYXRlbHkgYWZ0ZXIgdGhlc2UgcnVsaW5ncywgdGhlcmUgaXMgdGhlIG1pZ2h0IG9mIEFsbGFoLiAKaDM6IEJlbGlldmluZyBpbiBRdXJhbgpwOiBGaXJzdCBtaXJhY2xlIGluIGFueWJvZHnigJlzIGxpZmUgaXMsIG1heSBiZSwgYmlydGguIFlvdSBjb21lIG91dCBvZiB5b3VyIDxiPm1vdGhlcuKAmXMgYmVsbHk8L2I+LiBZb3UgZ3JvdyBhcyB5b3UgZWF0PGI+IGZvb2Q8L2I+LiBBZnRlciBzcGVuZGluZyB0aGUgcmVzdCBvZiB0aGUgbGlmZSBpbiB0aGlzIHdvcmxkLCB5b3U8Yj4gZGllPC9iPi4KcDogRHVyaW5nIGVhY2ggbWlyYWNsZSwgc2NpZW5jZSBzdHJpdmVzIHRvIGV4cGxhaW4gaXQuIFRoZSBiZWxseSBnaXZlcyB0aGUgY2hpbGQsIHdlIHN0YXJ0IGxhYmVsbGluZyBpdCBhIG5hdHVyYWwgYWN0LiBTY2llbmNlIGxhYmVscyBpdCBzZXF1ZW5jZSBvZiBtaXJhY3Vsb3VzIGV2ZW50cy4gCnA6IFdlIGxhYmVsIGZvb2QsIG5hdHVyYWwgaW5ncmVkaWVudHMsIHRvIHJ1biBodW1hbiBib2R5LiBXZSBncm93LCBsYWJlbGxpbmcgaXQgbmF0dXJhbC4gV2UgZGllIGFuZCBsYWJlbCBkZWF0aCBhIHJlc3VsdCBvZiBOYXR1cmFsIENhdXNlcy4KcDogVGhlIGVhcnRoIGlzIGFuIGV4dHJhIG9yZGluYXJ5IHBpZWNlIG9mIG5hdHVyZS4gRWFjaCBwYXNzaW5nIHNlY29uZCwgdGhlIGVhcnRoIGlzIHN1cnJvdW5kZWQgd2l0aCBhbm90aGVyIG1pcmFjbGUuIFdlIHdpdG5lc3MgdGhpcyBtaXJhY2xlIGFuZCA8Yj53ZSBydW4gdG8gcHJvdmUgaXQgaXMgbmF0dXJhbCBvciBzY2llbnRpZmljPC9iPi4gCnA6IFdlIHN0cml2ZSBoYXJkIHRvIG1ha2Ugb3Vyc2VsdmVzIGJlbGlldmUgaXQgaXMganVzdCBhbm90aGVyIGFjdCBvZiBuYXR1cmUuIApoMzogTmF0dXJlID0gPwpwOiBOYXR1cmUgaXMgbm90IGFydGlmaWNpYWwuIEl0IGlzIGFueXRoaW5nIGJyZWF0aGluZy4gVGhlIGVhcnRoLCBpcyBicmVhdGhpbmcuIEFuaW1hbHMgZ3JhemluZyBvbiBpdCBhcmUgYWxzbyBicmVhdGhpbmcuIFdlIGJlbGlldmUgdGhlIG5hdHVyZSBvZiBlYXJ0aCBpcyBjYXVzaW5nIHRoaXMgcGxhbmV0IHRvIHJvdGF0ZS4gVGhlcmUgaXMgZ3Jhdml0eSBsaW5rZWQgd2l0aCB0aGUgZWFydGgsIHRoYXQga2VlcHMgdGhlIGFuaW1hbHMgZnJvbSBmYWxsaW5nIG9mZiBpdC4gQW55dGhpbmcgdGhhdCBjYW4gbm90IGJlIHJlY29uc3RydWN0ZWQgaXMgbmF0dXJlLgpwOiBXZSBsYWJlbCBuYXR1cmUgYXMgYSByZXN1bHQgb2YgY29pbmNpZGVuY2Ugb2YgbWF0dGVyIGFuZCBtb2xlY3VsZS4KaDM6IERpZCBhbnlvbmUgZmluZCB0aGUgYW5zd2VyPwpwOiBFdmVyeSBjZWxsdWxhciBzdHJ1Y3R1cmUgaW5zaWRlIG5hdHVyZSBpcyBhIGNvbWJpbmF0aW9uIG9mIGVsZWN0cm9ucywgcHJvdG9ucyBhbmQgbmV1dHJvbnMuIFRoZXNlIGFyZSBtaW51dGUgcGFydGljbGVzIHRoYXQgbWFrZSBldmVyeSBtb2xlY3VsZS4gU29tZSBncm93IHRvIGJlY29tZSBhIHNwZXJtIGFuZCBmZXcgZ3JvdyB0byBiZWNvbWUgZWFydGguIEZyb20gbWludXRlIHRvIGJpZ2dlc3Qgc3RydWN0dXJlIGlzIGEgY29tYmluYXRpb24gb2YgY2VsbHMuIApwOiBOZXh0IHF1ZXN0aW9ucyB0aGF0IG5lZWQgdG8gYmUgYXNrZWQKcC5iLWxlZnQ6IDEuIFdoeSBzb21lIGVsZWN0cm9ucyBncm93IHRvIGxvb2sgbGlrZSBzcGVybT8KcC5iLWxlZnQ6IDIuIEhvdyBzcGVybSBkZWNpZGVzIGl0IGdvZXMgdG8gd29tYj8KcC5iLWxlZnQ6IDMuIEhvdyBkb2VzIGEgYmFieSBzdGF5cyBhbGl2ZSBpbnNpZGUgYSBjbG9zZWQgY2hhbWJlcj8KcC5iLWxlZnQ6IDQuIEhvdyBkb2VzIGEgc3Blcm0gZ3Jvd3MgdG8gYSBodW1hbiBiZWluZz8KcC5iLWxlZnQ6IDUuIFdoYXQga2VlcHMgdGhlIGh1bWFuIGZyb20gZmFsbGluZyBvZmYgdGhlIHBsYW5ldCBFYXJ0aD8KcC5iLWxlZnQ6IDYuIEhvdyBtYW55IG1vcmUgRWFydGjigJlzIGV4aXN0IGluIGEgdW5pdmVyc2Ugb2YgY29pbmNpZGVuY2VzPwpwLmItbGVmdDogNy4gSG93IGJpZyBpcyB1bml2ZXJzZT8KaDM6IERpZCB3ZSBmaW5kIHRoZSBhbnN3ZXIgdG8gb3VyIHF1ZXN0aW9ucz8KcDogTm8uIE91ciBubyBxdWVzdGlvbiBhYm91dCBuYXR1cmUgZ29lcyBhbnN3ZXJlZC4gRWFjaCBxdWVzdGlvbiB0YWtlcyB1cyBmdXJ0aGVyIGRvd24gYSBzZXJpZXMgb2YgcmFiYml0IGhvbGVzIGFuZCB0aGVzZSBxdWVzdGlvbnMgbmV2ZXIgc3RvcC4gIApwOiBXZSBkaWQgbm90IGZpbmQgb3V0IHRoZSByZWFsaXR5IGJlaGluZCBhIHNwZXJtLiBXaHkgaXQgZW50ZXJzIHRoZSB3b21iPyBXaGF0IG1ha2VzIGEgYmlydGg/IEVhY2ggcXVlc3Rpb24gaXMgb25seSBhbiBhcmd1bWVudC4KcDogVGhlIGFuc3dlcnMgd2UgZmluZCBhcmUgbGl0ZXJhbCBzY2llbnRpZmljIHN0b3JpZXMuIFNjaWVuY2UgbmV2ZXIgc3VycmVuZGVycywgaXQgaXMgYSBuYW1lIG9mIHJhYmJpdCBob2xlLiAKaDM6IENvbmNsdXNpb24KcDogV2UgZG8gbm90IGFzayAxIHF1ZXN0aW9uLiAKcC5iLWxlZnQ6IDxiPldoYXQgaWYgd2UgZGllIGFuZCB0aGVyZSBpcyBhIEdvZD88L2I+CnA6IFdlIGRvIG5vdCBzdWJtaXQgdG8gdGhlIHZhY3V1bSBjcmVhdGVkIGJ5IG91ciB1bmFuc3dlcmVkIHF1ZXN0aW9ucy4gRXhlY3V0aW9uIG9uIFF1cmFuaWMgaW5zdHJ1Y3Rpb25zIGlzIG9ubHkgcG9zc2libGUgYWZ0ZXIgd2UgZmluZCBBbGxhaC4gCnA6IE9uY2UgZm91bmQ7IHN1Ym1pdCwgc3VycmVuZGVyIGFuZCBtb3ZlIG5lYXIgSGltLiBIZSBpcyBnb2luZyB0byBtZWV0IHVzIHRoZSBkYXkgd2UgZGllIGluIGZldyBkYXlzIG9yIG1vbnRocy4gQSBNaWdodCBpcyBleHBlY3RpbmcgdXMuIDxiPkl0IGNhbiBoYXBwZW4gdGhhdCB0aGUgZGF5IG9mIG91ciBkZWF0aCBiZWNvbWVzIHRoZSDigJhIYXBwaWVzdCBkYXnigJkgb2Ygb3VyIGxpdmVzLiA8L2I+CnAubm90ZTogR3V5cywgaWYgeW91IGxpa2UgdGhpcyBwcm9qZWN0LiBQbGVhc2UgZm9sbG93IHRoaXMgcHJvamVjdCdzIHBhZ2Ugb24gdHdpdHRlci4gPGEgaHJlZj0iaHR0cHM6Ly90d2l0dGVyLmNvbS96YWthdGxpc3RzIj5DbGljayBoZXJlIHRvIGdvIHRvIHRoZSB0d2l0dGVyIHBhZ2U8L2E+LiBTdWJzY3JpYmUgYmVsb3cgdG8gcmVjZWl2ZSB1cGRhdGVzLmgxOiBNaWdodCBvZiBBbGxhaApwOiBUb2RheSBRdXJhbiB0YWxrcyBwcmVjaXNlbHkgYWJvdXQgdGhlIGRpc3RyaWJ1dGlvbiBvZiBpbmhlcml0YW5jZSB0byBvdXIgY2xvc2Ugb25lcy4gSWYgQ2hpbGQgZGllcywgbW90aGVyIGdldHMgMS82dGggb2YgaGVyIHNoYXJlLiBJZiBtb3RoZXIgaXMgYSB3aWRvd2VyLCBzaGUgZ2V0cyAxLzh0aCBvZiBoZXIgc2hhcmUuIEFmdGVyIGVhY2ggZGVjaXNpb24sIFF1cmFuIHNwZWFrcyBvZiB0aGUgbWlnaHQgb2YgdGhlIENyZWF0b3IuIEZvciBleGFtcGxlCnAuYi1sZWZ0OiA8Yj4xLiBWZXJzZSA0PC9iPi4gQW5kIHN1ZmZpY2llbnQgaXMgQWxsYWggYXMgQWNjb3VudGFudC4KcC5iLWxlZnQ6IDxiPjIuIFZlcnNlIDk8L2I+LiBTbyBsZXQgdGhlbSBmZWFyIEFsbGFoIGFuZCBzcGVhayB3b3JkcyBvZiBhcHByb3ByaWF0ZSBqdXN0aWNlLgpwLmItbGVmdDogPGI+My4gVmVyc2UgMTA8L2I+LiBBbmQgdGhleSB3aWxsIGJlIGJ1cm5lZCBpbiBhIEJsYXplLgpwLmItbGVmdDogPGI+NC4gVmVyc2UgMTE8L2I+LiBJbmRlZWQsIEFsbGFoIGlzIGV2ZXIgS25vd2luZyBhbmQgV2lzZS4KcC5iLWxlZnQ6IDxiPjUuIFZlcnNlIDEyPC9iPi4gQWxsYWggaXMgS25vd2luZyBhbmQgRm9yYmVhcmluZy4KcDogVGhlIHJoeXRobSBvZiB0aGVzZSBzZW50ZW5jZXMgaXMgZm91bmQgaW4gZW50aXJlIFF1cmFuLiBUaGVyZSBhcmUgZGVjaXNpb25zLCBpbnN0cnVjdGlvbnMsIG9ibGlnYXRpb25zIGFuZCBpbW1lZGlhdGVseSBhZnRlciB0aGVzZSBydWxpbmdzLCB0aGVyZSBpcyB0aGUgbWlnaHQgb2YgQWxsYWguIApoMzogQmVsaWV2aW5nIGluIFF1cmFuCnA6IEZpcnN0IG1pcmFjbGUgaW4gYW55Ym9keeKAmXMgbGlmZSBpcywgbWF5IGJlLCBiaXJ0aC4gWW91IGNvbWUgb3V0IG9mIHlvdXIgPGI+bW90aGVy4oCZcyBiZWxseTwvYj4uIFlvdSBncm93IGFzIHlvdSBlYXQ8Yj4gZm9vZDwvYj4uIEFmdGVyIHNwZW5kaW5nIHRoZSByZXN0IG9mIHRoZSBsaWZlIGluIHRoaXMgd29ybGQsIHlvdTxiPiBkaWU8L2I+LgpwOiBEdXJpbmcgZWFjaCBtaXJhY2xlLCBzY2llbmNlIHN0cml2ZXMgdG8gZXhwbGFpbiBpdC4gVGhlIGJlbGx5IGdpdmVzIHRoZSBjaGlsZCwgd2Ugc3RhcnQgbGFiZWxsaW5nIGl0IGEgbmF0dXJhbCBhY3QuIFNjaWVuY2UgbGFiZWxzIGl0IHNlcXVlbmNlIG9mIG1pcmFjdWxvdXMgZXZlbnRzLiAKcDogV2UgbGFiZWwgZm9vZCwgbmF0dXJhbCBpbmdyZWRpZW50cywgdG8gcnVuIGh1bWFuIGJvZHkuIFdlIGdyb3csIGxhYmVsbGluZyBpdCBuYXR1cmFsLiBXZSBkaWUgYW5kIGxhYmVsIGRlYXRoIGEgcmVzdWx0IG9mIE5hdHVyYWwgQ2F1c2VzLgpwOiBUaGUgZWFydGggaXMgYW4gZXh0cmEgb3JkaW5hcnkgcGllY2Ugb2YgbmF0dXJlLiBFYWNoIHBhc3Npbmcgc2Vjb25kLCB0aGUgZWFydGggaXMgc3Vycm91bmRlZCB3aXRoIGFub3RoZXIgbWlyYWNsZS4gV2Ugd2l0bmVzcyB0aGlzIG1pcmFjbGUgYW5kIDxiPndlIHJ1biB0byBwcm92ZSBpdCBpcyBuYXR1cmFsIG9yIHNjaWVudGlmaWM8L2I+LiAKcDogV2Ugc3RyaXZlIGhhcmQgdG8gbWFrZSBvdXJzZWx2ZXMgYmVsaWV2ZSBpdCBpcyBqdXN0IGFub3RoZXIgYWN0IG9mIG5hdHVyZS4gCmgzOiBOYXR1cmUgPSA/CnA6IE5hdHVyZSBpcyBub3QgYXJ0aWZpY2lhbC4gSXQgaXMgYW55dGhpbmcgYnJlYXRoaW5nLiBUaGUgZWFydGgsIGlzIGJyZWF0aGluZy4gQW5pbWFscyBncmF6aW5nIG9uIGl0IGFyZSBhbHNvIGJyZWF0aGluZy4gV2UgYmVsaWV2ZSB0aGUgbmF0dXJlIG9mIGVhcnRoIGlzIGNhdXNpbmcgdGhpcyBwbGFuZXQgdG8gcm90YXRlLiBUaGVyZSBpcyBncmF2aXR5IGxpbmtlZCB3aXRoIHRoZSBlYXJ0aCwgdGhhdCBrZWVwcyB0aGUgYW5pbWFscyBmcm9tIGZhbGxpbmcgb2ZmIGl0LiBBbnl0aGluZyB0aGF0IGNhbiBub3QgYmUgcmVjb25zdHJ1Y3RlZCBpcyBuYXR1cmUuCnA6IFdlIGxhYmVsIG5hdHVyZSBhcyBhIHJlc3VsdCBvZiBjb2luY2lkZW5jZSBvZiBtYXR0ZXIgYW5kIG1vbGVjdWxlLgpoMzogRGlkIGFueW9uZSBmaW5kIHRoZSBhbnN3ZXI/CnA6IEV2ZXJ5IGNlbGx1bGFyIHN0cnVjdHVyZSBpbnNpZGUgbmF0dXJlIGlzIGEgY29tYmluYXRpb24gb2YgZWxlY3Ryb25zLCBwcm90b25zIGFuZCBuZXV0cm9ucy4gVGhlc2UgYXJlIG1pbnV0ZSBwYXJ0aWNsZXMgdGhhdCBtYWtlIGV2ZXJ5IG1vbGVjdWxlLiBTb21lIGdyb3cgdG8gYmVjb21lIGEgc3Blcm0gYW5kIGZldyBncm93IHRvIGJlY29tZSBlYXJ0aC4gRnJvbSBtaW51dGUgdG8gYmlnZ2VzdCBzdHJ1Y3R1cmUgaXMgYSBjb21iaW5hdGlvbiBvZiBjZWxscy4gCnA6IE5leHQgcXVlc3Rpb25zIHRoYXQgbmVlZCB0byBiZSBhc2tlZApwLmItbGVmdDogMS4gV2h5IHNvbWUgZWxlY3Ryb25zIGdyb3cgdG8gbG9vayBsaWtlIHNwZXJtPwpwLmItbGVmdDogMi4gSG93IHNwZXJtIGRlY2lkZXMgaXQgZ29lcyB0byB3b21iPwpwLmItbGVmdDogMy4gSG93IGRvZXMgYSBiYWJ5IHN0YXlzIGFsaXZlIGluc2lkZSBhIGNsb3NlZCBjaGFtYmVyPwpwLmItbGVmdDogNC4gSG93IGRvZXMgYSBzcGVybSBncm93cyB0byBhIGh1bWFuIGJlaW5nPwpwLmItbGVmdDogNS4gV2hhdCBrZWVwcyB0aGUgaHVtYW4gZnJvbSBmYWxsaW5nIG9mZiB0aGUgcGxhbmV0IEVhcnRoPwpwLmItbGVmdDogNi4gSG93IG1hbnkgbW9yZSBFYXJ0aOKAmXMgZXhpc3QgaW4gYSB1bml2ZXJzZSBvZiBjb2luY2lkZW5jZXM/CnAuYi1sZWZ0OiA3LiBIb3cgYmlnIGlzIHVuaXZlcnNlPwpoMzogRGlkIHdlIGZpbmQgdGhlIGFuc3dlciB0byBvdXIgcXVlc3Rpb25zPwpwOiBOby4gT3VyIG5vIHF1ZXN0aW9uIGFib3V0IG5hdHVyZSBnb2VzIGFuc3dlcmVkLiBFYWNoIHF1ZXN0aW9uIHRha2VzIHVzIGZ1cnRoZXIgZG93biBhIHNlcmllcyBvZiByYWJiaXQgaG9sZXMgYW5kIHRoZXNlIHF1ZXN0aW9ucyBuZXZlciBzdG9wLiAgCnA6IFdlIGRpZCBub3QgZmluZCBvdXQgdGhlIHJlYWxpdHkgYmVoaW5kIGEgc3Blcm0uIFdoeSBpdCBlbnRlcnMgdGhlIHdvbWI/IFdoYXQgbWFrZXMgYSBiaXJ0aD8gRWFjaCBxdWVzdGlvbiBpcyBvbmx5IGFuIGFyZ3VtZW50LgpwOiBUaGUgYW5zd2VycyB3ZSBmaW5kIGFyZSBsaXRlcmFsIHNjaWVudGlmaWMgc3Rvcmllcy4gU2NpZW5jZSBuZXZlciBzdXJyZW5kZXJzLCBpdCBpcyBhIG5hbWUgb2YgcmFiYml0IGhvbGUuIApoMzogQ29uY2x1c2lvbgpwOiBXZSBkbyBub3QgYXNrIDEgcXVlc3Rpb24uIApwLmItbGVmdDogPGI+V2hhdCBpZiB3ZSBkaWUgYW5kIHRoZXJlIGlzIGEgR29kPzwvYj4KcDogV2UgZG8gbm90IHN1Ym1pdCB0byB0aGUgdmFjdXVtIGNyZWF0ZWQgYnkgb3VyIHVuYW5zd2VyZWQgcXVlc3Rpb25zLiBFeGVjdXRpb24gb24gUXVyYW5pYyBpbnN0cnVjdGlvbnMgaXMgb25seSBwb3NzaWJsZSBhZnRlciB3ZSBmaW5kIEFsbGFoLiAKcDogT25jZSBmb3VuZDsgc3VibWl0LCBzdXJyZW5kZXIgYW5kIG1vdmUgbmVhciBIaW0uIEhlIGlzIGdvaW5nIHRvIG1lZXQgdXMgdGhlIGRheSB3ZSBkaWUgaW4gZmV3IGRheXMgb3IgbW9udGhzLiBBIE1pZ2h0IGlzIGV4cGVjdGluZyB1cy4gPGI+SXQgY2FuIGhhcHBlbiB0aGF0IHRoZSBkYXkgb2Ygb3VyIGRlYXRoIGJlY29tZXMgdGhlIOKAmEhhcHBpZXN0IGRheeKAmSBvZiBvdXIgbGl2ZXMuIDwvYj4KcC5ub3RlOiBHdXlzLCBpZiB5b3UgbGlrZSB0aGlzIHByb2plY3QuIFBsZWFzZSBmb2xsb3cgdGhpcyBwcm9qZWN0J3MgcGFnZSBvbiB0d2l0dGVyLiA8YSBocmVmPSJodHRwczovL3R3aXR0ZXIuY29tL3pha2F0bGlzdHMiPkNsaWNrIGhlcmUgdG8gZ28gdG8gdGhlIHR3aXR0ZXIgcGFnZTwvYT4uIFN1YnNjcmliZSBiZWxvdyB0byByZWNlaXZlIHVwZGF0ZXMuaDE6IE1pZ2h0IG9mIEFsbGFoCnA6IFRvZGF5IFF1cmFuIHRhbGtzIHByZWNpc2VseSBhYm91dCB0aGUgZGlzdHJpYnV0aW9uIG9mIGluaGVyaXRhbmNlIHRvIG91ciBjbG9zZSBvbmVzLiBJZiBDaGlsZCBkaWVzLCBtb3RoZXIgZ2V0cyAxLzZ0aCBvZiBoZXIgc2hhcmUuIElmIG1vdGhlciBpcyBhIHdpZG93ZXIsIHNoZSBnZXRzIDEvOHRoIG9mIGhlciBzaGFyZS4gQWZ0ZXIgZWFjaCBkZWNpc2lvbiwgUXVyYW4gc3BlYWtzIG9mIHRoZSBtaWdodCBvZiB0aGUgQ3JlYXRvci4gRm9yIGV4YW1wbGUKcC5iLWxlZnQ6IDxiPjEuIFZlcnNlIDQ8L2I+LiBBbmQgc3VmZmljaWVudCBpcyBBbGxhaCBhcyBBY2NvdW50YW50LgpwLmItbGVmdDogPGI+Mi4gVmVyc2UgOTwvYj4uIFNvIGxldCB0aGVtIGZlYXIgQWxsYWggYW5kIHNwZWFrIHdvcmRzIG9mIGFwcHJvcHJpYXRlIGp1c3RpY2UuCnAuYi1sZWZ0OiA8Yj4zLiBWZXJzZSAxMDwvYj4uIEFuZCB0aGV5IHdpbGwgYmUgYnVybmVkIGluIGEgQmxhemUuCnAuYi1sZWZ0OiA8Yj40LiBWZXJzZSAxMTwvYj4uIEluZGVlZCwgQWxsYWggaXMgZXZlciBLbm93aW5nIGFuZCBXaXNlLgpwLmItbGVmdDogPGI+NS4gVmVyc2UgMTI8L2I+LiBBbGxhaCBpcyBLbm93aW5nIGFuZCBGb3JiZWFyaW5nLgpwOiBUaGUgcmh5dGhtIG9mIHRoZXNlIHNlbnRlbmNlcyBpcyBmb3VuZCBpbiBlbnRpcmUgUXVyYW4uIFRoZXJlIGFyZSBkZWNpc2lvbnMsIGluc3RydWN0aW9ucywgb2JsaWdhdGlvbnMgYW5kIGltbWVkaWF0ZWx5IGFmdGVyIHRoZXNlIHJ1bGluZ3MsIHRoZXJlIGlzIHRoZSBtaWdodCBvZiBBbGxhaC4gCmgzOiBCZWxpZXZpbmcgaW4gUXVyYW4KcDogRmlyc3QgbWlyYWNsZSBpbiBhbnlib2R54oCZcyBsaWZlIGlzLCBtYXkgYmUsIGJpcnRoLiBZb3UgY29tZSBvdXQgb2YgeW91ciA8Yj5tb3RoZXLigJlzIGJlbGx5PC9iPi4gWW91IGdyb3cgYXMgeW91IGVhdDxiPiBmb29kPC9iPi4gQWZ0ZXIgc3BlbmRpbmcgdGhlIHJlc3Qgb2YgdGhlIGxpZmUgaW4gdGhpcyB3b3JsZCwgeW91PGI+IGRpZTwvYj4uCnA6IER1cmluZyBlYWNoIG1pcmFjbGUsIHNjaWVuY2Ugc3RyaXZlcyB0byBleHBsYWluIGl0LiBUaGUgYmVsbHkgZ2l2ZXMgdGhlIGNoaWxkLCB3ZSBzdGFydCBsYWJlbGxpbmcgaXQgYSBuYXR1cmFsIGFjdC4gU2NpZW5jZSBsYWJlbHMgaXQgc2VxdWVuY2Ugb2YgbWlyYWN1bG91cyBldmVudHMuIApwOiBXZSBsYWJlbCBmb29kLCBuYXR1cmFsIGluZ3JlZGllbnRzLCB0byBydW4gaHVtYW4gYm9keS4gV2UgZ3JvdywgbGFiZWxsaW5nIGl0IG5hdHVyYWwuIFdlIGRpZSBhbmQgbGFiZWwgZGVhdGggYSByZXN1bHQgb2YgTmF0dXJhbCBDYXVzZXMuCnA6IFRoZSBlYXJ0aCBpcyBhbiBleHRyYSBvcmRpbmFyeSBwaWVjZSBvZiBuYXR1cmUuIEVhY2ggcGFzc2luZyBzZWNvbmQsIHRoZSBlYXJ0aCBpcyBzdXJyb3VuZGVkIHdpdGggYW5vdGhlciBtaXJhY2xlLiBXZSB3aXRuZXNzIHRoaXMgbWlyYWNsZSBhbmQgPGI+d2UgcnVuIHRvIHByb3ZlIGl0IGlzIG5hdHVyYWwgb3Igc2NpZW50aWZpYzwvYj4uIApwOiBXZSBzdHJpdmUgaGFyZCB0byBtYWtlIG91cnNlbHZlcyBiZWxpZXZlIGl0IGlzIGp1c3QgYW5vdGhlciBhY3Qgb2YgbmF0dXJlLiAKaDM6IE5hdHVyZSA9ID8KcDogTmF0dXJlIGlzIG5vdCBhcnRpZmljaWFsLiBJdCBpcyBhbnl0aGluZyBicmVhdGhpbmcuIFRoZSBlYXJ0aCwgaXMgYnJlYXRoaW5nLiBBbmltYWxzIGdyYXppbmcgb24gaXQgYXJlIGFsc28gYnJlYXRoaW5nLiBXZSBiZWxpZXZlIHRoZSBuYXR1cmUgb2YgZWFydGggaXMgY2F1c2luZyB0aGlzIHBsYW5ldCB0byByb3RhdGUuIFRoZXJlIGlzIGdyYXZpdHkgbGlua2VkIHdpdGggdGhlIGVhcnRoLCB0aGF0IGtlZXBzIHRoZSBhbmltYWxzIGZyb20gZmFsbGluZyBvZmYgaXQuIEFueXRoaW5nIHRoYXQgY2FuIG5vdCBiZSByZWNvbnN0cnVjdGVkIGlzIG5hdHVyZS4KcDogV2UgbGFiZWwgbmF0dXJlIGFzIGEgcmVzdWx0IG9mIGNvaW5jaWRlbmNlIG9mIG1hdHRlciBhbmQgbW9sZWN1bGUuCmgzOiBEaWQgYW55b25lIGZpbmQgdGhlIGFuc3dlcj8KcDogRXZlcnkgY2VsbHVsYXIgc3RydWN0dXJlIGluc2lkZSBuYXR1cmUgaXMgYSBjb21iaW5hdGlvbiBvZiBlbGVjdHJvbnMsIHByb3RvbnMgYW5kIG5ldXRyb25zLiBUaGVzZSBhcmUgbWludXRlIHBhcnRpY2xlcyB0aGF0IG1ha2UgZXZlcnkgbW9sZWN1bGUuIFNvbWUgZ3JvdyB0byBiZWNvbWUgYSBzcGVybSBhbmQgZmV3IGdyb3cgdG8gYmVjb21lIGVhcnRoLiBGcm9tIG1pbnV0ZSB0byBiaWdnZXN0IHN0cnVjdHVyZSBpcyBhIGNvbWJpbmF0aW9uIG9mIGNlbGxzLiAKcDogTmV4dCBxdWVzdGlvbnMgdGhhdCBuZWVkIHRvIGJlIGFza2VkCnAuYi1sZWZ0OiAxLiBXaHkgc29tZSBlbGVjdHJvbnMgZ3JvdyB0byBsb29rIGxpa2Ugc3Blcm0/CnAuYi1sZWZ0OiAyLiBIb3cgc3Blcm0gZGVjaWRlcyBpdCBnb2VzIHRvIHdvbWI/CnAuYi1sZWZ0OiAzLiBIb3cgZG9lcyBhIGJhYnkgc3RheXMgYWxpdmUgaW5zaWRlIGEgY2xvc2VkIGNoYW1iZXI/CnAuYi1sZWZ0OiA0LiBIb3cgZG9lcyBhIHNwZXJtIGdyb3dzIHRvIGEgaHVtYW4gYmVpbmc/CnAuYi1sZWZ0OiA1LiBXaGF0IGtlZXBzIHRoZSBodW1hbiBmcm9tIGZhbGxpbmcgb2ZmIHRoZSBwbGFuZXQgRWFydGg/CnAuYi1sZWZ0OiA2LiBIb3cgbWFueSBtb3JlIEVhcnRo4oCZcyBleGlzdCBpbiBhIHVuaXZlcnNlIG9mIGNvaW5jaWRlbmNlcz8KcC5iLWxlZnQ6IDcuIEhvdyBiaWcgaXMgdW5pdmVyc2U/CmgzOiBEaWQgd2UgZmluZCB0aGUgYW5zd2VyIHRvIG91ciBxdWVzdGlvbnM/CnA6IE5vLiBPdXIgbm8gcXVlc3Rpb24gYWJvdXQgbmF0dXJlIGdvZXMgYW5zd2VyZWQuIEVhY2ggcXVlc3Rpb24gdGFrZXMgdXMgZnVydGhlciBkb3duIGEgc2VyaWVzIG9mIHJhYmJpdCBob2xlcyBhbmQgdGhlc2UgcXVlc3Rpb25zIG5ldmVyIHN0b3AuICAKcDogV2UgZGlkIG5vdCBmaW5kIG91dCB0aGUgcmVhbGl0eSBiZWhpbmQgYSBzcGVybS4gV2h5IGl0IGVudGVycyB0aGUgd29tYj8gV2hhdCBtYWtlcyBhIGJpcnRoPyBFYWNoIHF1ZXN0aW9uIGlzIG9ubHkgYW4gYXJndW1lbnQuCnA6IFRoZSBhbnN3ZXJzIHdlIGZpbmQgYXJlIGxpdGVyYWwgc2NpZW50aWZpYyBzdG9yaWVzLiBTY2llbmNlIG5ldmVyIHN1cnJlbmRlcnMsIGl0IGlzIGEgbmFtZSBvZiByYWJiaXQgaG9sZS4gCmgzOiBDb25jbHVzaW9uCnA6IFdlIGRvIG5vdCBhc2sgMSBxdWVzdGlvbi4gCnAuYi1sZWZ0OiA8Yj5XaGF0IGlmIHdlIGRpZSBhbmQgdGhlcmUgaXMgYSBHb2Q/PC9iPgpwOiBXZSBkbyBub3Qgc3VibWl0IHRvIHRoZSB2YWN1dW0gY3JlYXRlZCBieSBvdXIgdW5hbnN3ZXJlZCBxdWVzdGlvbnMuIEV4ZWN1dGlvbiBvbiBRdXJhbmljIGluc3RydWN0aW9ucyBpcyBvbmx5IHBvc3NpYmxlIGFmdGVyIHdlIGZpbmQgQWxsYWguIApwOiBPbmNlIGZvdW5kOyBzdWJtaXQsIHN1cnJlbmRlciBhbmQgbW92ZSBuZWFyIEhpbS4gSGUgaXMgZ29pbmcgdG8gbWVldCB1cyB0aGUgZGF5IHdlIGRpZSBpbiBmZXcgZGF5cyBvciBtb250aHMuIEEgTWlnaHQgaXMgZXhwZWN0aW5nIHVzLiA8Yj5JdCBjYW4gaGFwcGVuIHRoYXQgdGhlIGRheSBvZiBvdXIgZGVhdGggYmVjb21lcyB0aGUg4oCYSGFwcGllc3QgZGF54oCZIG9mIG91ciBsaXZlcy4gPC9iPgpwLm5vdGU6IEd1eXMsIGlmIHlvdSBsaWtlIHRoaXMgcHJvamVjdC4gUGxlYXNlIGZvbGxvdyB0aGlzIHByb2plY3QncyBwYWdlIG9uIHR3aXR0ZXIuIDxhIGhyZWY9Imh0dHBzOi8vdHdpdHRlci5jb20vemFrYXRsaXN0cyI+Q2xpY2sgaGVyZSB0byBnbyB0byB0aGUgdHdpdHRlciBwYWdlPC9hPi4gU3Vic2NyaWJlIGJlbG93IHRvIHJlY2VpdmUgdXBkYXRlcy5oMTogTWlnaHQgb2YgQWxsYWgKcDogVG9kYXkgUXVyYW4gdGFsa3MgcHJlY2lzZWx5IGFib3V0IHRoZSBkaXN0cmlidXRpb24gb2YgaW5oZXJpdGFuY2UgdG8gb3VyIGNsb3NlIG9uZXMuIElmIENoaWxkIGRpZXMsIG1vdGhlciBnZXRzIDEvNnRoIG9mIGhlciBzaGFyZS4gSWYgbW90aGVyIGlzIGEgd2lkb3dlciwgc2hlIGdldHMgMS84dGggb2YgaGVyIHNoYXJlLiBBZnRlciBlYWNoIGRlY2lzaW9uLCBRdXJhbiBzcGVha3Mgb2YgdGhlIG1pZ2h0IG9mIHRoZSBDcmVhdG9yLiBGb3IgZXhhbXBsZQpwLmItbGVmdDogPGI+MS4gVmVyc2UgNDwvYj4uIEFuZCBzdWZmaWNpZW50IGlzIEFsbGFoIGFzIEFjY291bnRhbnQuCnAuYi1sZWZ0OiA8Yj4yLiBWZXJzZSA5PC9iPi4gU28gbGV0IHRoZW0gZmVhciBBbGxhaCBhbmQgc3BlYWsgd29yZHMgb2YgYXBwcm9wcmlhdGUganVzdGljZS4KcC5iLWxlZnQ6IDxiPjMuIFZlcnNlIDEwPC9iPi4gQW5kIHRoZXkgd2lsbCBiZSBidXJuZWQgaW4gYSBCbGF6ZS4KcC5iLWxlZnQ6IDxiPjQuIFZlcnNlIDExPC9iPi4gSW5kZWVkLCBBbGxhaCBpcyBldmVyIEtub3dpbmcgYW5kIFdpc2UuCnAuYi1sZWZ0OiA8Yj41LiBWZXJzZSAxMjwvYj4uIEFsbGFoIGlzIEtub3dpbmcgYW5kIEZvcmJlYXJpbmcuCnA6IFRoZSByaHl0aG0gb2YgdGhlc2Ugc2VudGVuY2VzIGlzIGZvdW5kIGluIGVudGlyZSBRdXJhbi4gVGhlcmUgYXJlIGRlY2lzaW9ucywgaW5zdHJ1Y3Rpb25zLCBvYmxpZ2F0aW9ucyBhbmQgaW1tZWRpYXRlbHkgYWZ0ZXIgdGhlc2UgcnVsaW5ncywgdGhlcmUgaXMgdGhlIG1pZ2h0IG9mIEFsbGFoLiAKaDM6IEJlbGlldmluZyBpbiBRdXJhbgpwOiBGaXJzdCBtaXJhY2xlIGluIGFueWJvZHnigJlzIGxpZmUgaXMsIG1heSBiZSwgYmlydGguIFlvdSBjb21lIG91dCBvZiB5b3VyIDxiPm1vdGhlcuKAmXMgYmVsbHk8L2I+LiBZb3UgZ3JvdyBhcyB5b3UgZWF0PGI+IGZvb2Q8L2I+LiBBZnRlciBzcGVuZGluZyB0aGUgcmVzdCBvZiB0aGUgbGlmZSBpbiB0aGlzIHdvcmxkLCB5b3U8Yj4gZGllPC9iPi4KcDogRHVyaW5nIGVhY2ggbWlyYWNsZSwgc2NpZW5jZSBzdHJpdmVzIHRvIGV4cGxhaW4gaXQuIFRoZSBiZWxseSBnaXZlcyB0aGUgY2hpbGQsIHdlIHN0YXJ0IGxhYmVsbGluZyBpdCBhIG5hdHVyYWwgYWN0LiBTY2llbmNlIGxhYmVscyBpdCBzZXF1ZW5jZSBvZiBtaXJhY3Vsb3VzIGV2ZW50cy4gCnA6IFdlIGxhYmVsIGZvb2QsIG5hdHVyYWwgaW5ncmVkaWVudHMsIHRvIHJ1biBodW1hbiBib2R5LiBXZSBncm93LCBsYWJlbGxpbmcgaXQgbmF0dXJhbC4gV2UgZGllIGFuZCBsYWJlbCBkZWF0aCBhIHJlc3VsdCBvZiBOYXR1cmFsIENhdXNlcy4KcDogVGhlIGVhcnRoIGlzIGFuIGV4dHJhIG9yZGluYXJ5IHBpZWNlIG9mIG5hdHVyZS4gRWFjaCBwYXNzaW5nIHNlY29uZCwgdGhlIGVhcnRoIGlzIHN1cnJvdW5kZWQgd2l0aCBhbm90aGVyIG1pcmFjbGUuIFdlIHdpdG5lc3MgdGhpcyBtaXJhY2xlIGFuZCA8Yj53ZSBydW4gdG8gcHJvdmUgaXQgaXMgbmF0dXJhbCBvciBzY2llbnRpZmljPC9iPi4gCnA6IFdlIHN0cml2ZSBoYXJkIHRvIG1ha2Ugb3Vyc2VsdmVzIGJlbGlldmUgaXQgaXMganVzdCBhbm90aGVyIGFjdCBvZiBuYXR1cmUuIApoMzogTmF0dXJlID0gPwpwOiBOYXR1cmUgaXMgbm90IGFydGlmaWNpYWwuIEl0IGlzIGFueXRoaW5nIGJyZWF0aGluZy4gVGhlIGVhcnRoLCBpcyBicmVhdGhpbmcuIEFuaW1hbHMgZ3JhemluZyBvbiBpdCBhcmUgYWxzbyBicmVhdGhpbmcuIFdlIGJlbGlldmUgdGhlIG5hdHVyZSBvZiBlYXJ0aCBpcyBjYXVzaW5nIHRoaXMgcGxhbmV0IHRvIHJvdGF0ZS4gVGhlcmUgaXMgZ3Jhdml0eSBsaW5rZWQgd2l0aCB0aGUgZWFydGgsIHRoYXQga2VlcHMgdGhlIGFuaW1hbHMgZnJvbSBmYWxsaW5nIG9mZiBpdC4gQW55dGhpbmcgdGhhdCBjYW4gbm90IGJlIHJlY29uc3RydWN0ZWQgaXMgbmF0dXJlLgpwOiBXZSBsYWJlbCBuYXR1cmUgYXMgYSByZXN1bHQgb2YgY29pbmNpZGVuY2Ugb2YgbWF0dGVyIGFuZCBtb2xlY3VsZS4KaDM6IERpZCBhbnlvbmUgZmluZCB0aGUgYW5zd2VyPwpwOiBFdmVyeSBjZWxsdWxhciBzdHJ1Y3R1cmUgaW5zaWRlIG5hdHVyZSBpcyBhIGNvbWJpbmF0aW9uIG9mIGVsZWN0cm9ucywgcHJvdG9ucyBhbmQgbmV1dHJvbnMuIFRoZXNlIGFyZSBtaW51dGUgcGFydGljbGVzIHRoYXQgbWFrZSBldmVyeSBtb2xlY3VsZS4gU29tZSBncm93IHRvIGJlY29tZSBhIHNwZXJtIGFuZCBmZXcgZ3JvdyB0byBiZWNvbWUgZWFydGguIEZyb20gbWludXRlIHRvIGJpZ2dlc3Qgc3RydWN0dXJlIGlzIGEgY29tYmluYXRpb24gb2YgY2VsbHMuIApwOiBOZXh0IHF1ZXN0aW9ucyB0aGF0IG5lZWQgdG8gYmUgYXNrZWQKcC5iLWxlZnQ6IDEuIFdoeSBzb21lIGVsZWN0cm9ucyBncm93IHRvIGxvb2sgbGlrZSBzcGVybT8KcC5iLWxlZnQ6IDIuIEhvdyBzcGVybSBkZWNpZGVzIGl0IGdvZXMgdG8gd29tYj8KcC5iLWxlZnQ6IDMuIEhvdyBkb2VzIGEgYmFieSBzdGF5cyBhbGl2ZSBpbnNpZGUgYSBjbG9zZWQgY2hhbWJlcj8KcC5iLWxlZnQ6IDQuIEhvdyBkb2VzIGEgc3Blcm0gZ3Jvd3MgdG8gYSBodW1hbiBiZWluZz8KcC5iLWxlZnQ6IDUuIFdoYXQga2VlcHMgdGhlIGh1bWFuIGZyb20gZmFsbGluZyBvZmYgdGhlIHBsYW5ldCBFYXJ0aD8KcC5iLWxlZnQ6IDYuIEhvdyBtYW55IG1vcmUgRWFydGjigJlzIGV4aXN0IGluIGEgdW5pdmVyc2Ugb2YgY29pbmNpZGVuY2VzPwpwLmItbGVmdDogNy4gSG93IGJpZyBpcyB1bml2ZXJzZT8KaDM6IERpZCB3ZSBmaW5kIHRoZSBhbnN3ZXIgdG8gb3VyIHF1ZXN0aW9ucz8KcDogTm8uIE91ciBubyBxdWVzdGlvbiBhYm91dCBuYXR1cmUgZ29lcyBhbnN3ZXJlZC4gRWFjaCBxdWVzdGlvbiB0YWtlcyB1cyBmdXJ0aGVyIGRvd24gYSBzZXJpZXMgb2YgcmFiYml0IGhvbGVzIGFuZCB0aGVzZSBxdWVzdGlvbnMgbmV2ZXIgc3RvcC4gIApwOiBXZSBkaWQgbm90IGZpbmQgb3V0IHRoZSByZWFsaXR5IGJlaGluZCBhIHNwZXJtLiBXaHkgaXQgZW50ZXJzIHRoZSB3b21iPyBXaGF0IG1ha2VzIGEgYmlydGg/IEVhY2ggcXVlc3Rpb24gaXMgb25seSBhbiBhcmd1bWVudC4KcDogVGhlIGFuc3dlcnMgd2UgZmluZCBhcmUgbGl0ZXJhbCBzY2llbnRpZmljIHN0b3JpZXMuIFNjaWVuY2UgbmV2ZXIgc3VycmVuZGVycywgaXQgaXMgYSBuYW1lIG9mIHJhYmJpdCBob2xlLiAKaDM6IENvbmNsdXNpb24KcDogV2UgZG8gbm90IGFzayAxIHF1ZXN0aW9uLiAKcC5iLWxlZnQ6IDxiPldoYXQgaWYgd2UgZGllIGFuZCB0aGVyZSBpcyBhIEdvZD88L2I+CnA6IFdlIGRvIG5vdCBzdWJtaXQgdG8gdGhlIHZhY3V1bSBjcmVhdGVkIGJ5IG91ciB1bmFuc3dlcmVkIHF1ZXN0aW9ucy4gRXhlY3V0aW9uIG9uIFF1cmFuaWMgaW5zdHJ1Y3Rpb25zIGlzIG9ubHkgcG9zc2libGUgYWZ0ZXIgd2UgZmluZCBBbGxhaC4gCnA6IE9uY2UgZm91bmQ7IHN1Ym1pdCwgc3VycmVuZGVyIGFuZCBtb3ZlIG5lYXIgSGltLiBIZSBpcyBnb2luZyB0byBtZWV0IHVzIHRoZSBkYXkgd2UgZGllIGluIGZldyBkYXlzIG9yIG1vbnRocy4gQSBNaWdodCBpcyBleHBlY3RpbmcgdXMuIDxiPkl0IGNhbiBoYXBwZW4gdGhhdCB0aGUgZGF5IG9mIG91ciBkZWF0aCBiZWNvbWVzIHRoZSDigJhIYXBwaWVzdCBkYXnigJkgb2Ygb3VyIGxpdmVzLiA8L2I+CnAubm90ZTogR3V5cywgaWYgeW91IGxpa2UgdGhpcyBwcm9qZWN0LiBQbGVhc2UgZm9sbG93IHRoaXMgcHJvamVjdCdzIHBhZ2Ugb24gdHdpdHRlci4gPGEgaHJlZj0iaHR0cHM6Ly90d2l0dGVyLmNvbS96YWthdGxpc3RzIj5DbGljayBoZXJlIHRvIGdvIHRvIHRoZSB0d2l0dGVyIHBhZ2U8L2E+LiBTdWJzY3JpYmUgYmVsb3cgdG8gcmVjZWl2ZSB1cGRhdGVzLg=="
$bytes = [System.Convert]::FromBase64String($b64)
$text = [System.Text.Encoding]::UTF8.GetString($bytes)
$ws.Range("D43").Value = $text
